# Apply weekly Fruta/Hortaliza update: insert 5 new price rows (Comercializadora
# del Agro de Limarí - Palta, Edranol Especial/Primera/Segunda + Hass 1a/2a nueva(o))
# right before the existing row 699, shifting all subsequent rows down by 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows at 699, pushing old 699..746 down to 704..751.
$ws.Rows("699:703").Insert()

# Shared/constant values for this market+product block.
$mercadoId = 2
$mercado   = "Comercializadora del Agro de Limarí"
$region    = "Coquimbo"
$codreg    = 4
$tipo      = "Fruta"
$productoId = 100106
$producto  = "Oleaginosos"
$categoriaId = 100106002
$categoria = "Palta"
$unidad    = "$/kilo (en caja de 17 kilos)"
$origen    = "Provincia de Limarí"
$kgUnidad  = 1

# New rows data: Row, Fecha, Variedad, Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm
$newRows = @(
    @(699, 45147, "Edranol", "Especial",    200, 2300, 2400, 2350),
    @(700, 45147, "Edranol", "Primera",     240, 2000, 2100, 2050),
    @(701, 45147, "Edranol", "Segunda",     160, 1700, 1800, 1750),
    @(702, 45147, "Hass",    "1a nueva(o)", 200, 2000, 2100, 2050),
    @(703, 45147, "Hass",    "2a nueva(o)", 160, 1700, 1800, 1750)
)

foreach ($r in $newRows) {
    $row = $r[0]
    $fecha = $r[1]
    $variedad = $r[2]
    $calidad = $r[3]
    $volumen = $r[4]
    $precioMin = $r[5]
    $precioMax = $r[6]
    $precioProm = $r[7]

    $ws.Cells.Item($row, 1).Value2  = $mercadoId
    $ws.Cells.Item($row, 2).Value2  = $mercado
    $ws.Cells.Item($row, 3).Value2  = $region
    $ws.Cells.Item($row, 4).Value2  = $fecha
    $ws.Cells.Item($row, 5).Value2  = $codreg
    $ws.Cells.Item($row, 6).Value2  = $tipo
    $ws.Cells.Item($row, 7).Value2  = $productoId
    $ws.Cells.Item($row, 8).Value2  = $producto
    $ws.Cells.Item($row, 9).Value2  = $categoriaId
    $ws.Cells.Item($row, 10).Value2 = $categoria
    $ws.Cells.Item($row, 11).Value2 = $variedad
    $ws.Cells.Item($row, 12).Value2 = $calidad
    $ws.Cells.Item($row, 13).Value2 = $volumen
    $ws.Cells.Item($row, 14).Value2 = $precioMin
    $ws.Cells.Item($row, 15).Value2 = $precioMax
    $ws.Cells.Item($row, 16).Value2 = $precioProm
    $ws.Cells.Item($row, 17).Value2 = $unidad
    $ws.Cells.Item($row, 18).Value2 = $origen
    $ws.Cells.Item($row, 19).Value2 = $precioProm
    $ws.Cells.Item($row, 20).Value2 = $kgUnidad
}
